$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G2").Value = 'Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Nahla Nagiub'
$ws.Range("G3").Value = 'Dr. Asmaa Reda, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid'
$ws.Range("G4").Value = 'Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Manar Montaser, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Majorelle Magdy'
$ws.Range("G5").Value = 'Dr. Lamiaa Ossama, Dr. Nada Gouda, Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady'
$ws.Range("G6").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range("G7").Value = 'Dr. Arwa Elnagar, Dr. Shimaa Ashraf, Dr. Aya Saeed'
$ws.Range("G8").Value = 'Dr. Dalia Mohammad Abd Al-Salam, Dr. Marwa Mustafa, Dr. Madeha Saeed, Dr. Dina Adel, Dr. Amira Ibrahim'
$ws.Range("G9").Value = 'Dr. Maryam Ahmad, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Merna Said'
$ws.Range("G10").Value = 'Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat'
$ws.Range("G15").Value = 'Dr. Walaa Ghanima, Dr. Amr Saeed'
$ws.Range("G16").Value = 'Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Aya Hanafy, Dr. Wafaa Ebida'
$ws.Range("G17").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G18").Value = 'Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Nahla Nagiub'
$ws.Range("G19").Value = 'Dr. Mohammad El-Tanany, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad'
$ws.Range("G20").Value = 'Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G21").Value = 'Dr. Lamiaa Ossama, Dr. Nada Gouda, Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady'
$ws.Range("G22").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range("G23").Value = 'Dr. Arwa Elnagar, Dr. Shimaa Ashraf, Dr. Aya Saeed'
$ws.Range("G24").Value = 'Dr. Dalia Mohammad Abd Al-Salam, Dr. Marwa Mustafa, Dr. Madeha Saeed, Dr. Dina Adel, Dr. Amira Ibrahim'
$ws.Range("G25").Value = 'Dr. Maryam Ahmad, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Merna Said'
$ws.Range("G26").Value = 'Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat'
$ws.Range("G31").Value = 'Dr. Walaa Ghanima, Dr. Amr Saeed'
$ws.Range("G32").Value = 'Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Aya Hanafy, Dr. Wafaa Ebida'
$ws.Range("G33").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G34").Value = 'Dr. Amira Sobhy, Administrator, Dr. Asmaa Reda, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Nahla Nagiub'
$ws.Range("G35").Value = 'Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Administrator, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Rana Abo-Zaid'
$ws.Range("G36").Value = 'Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G37").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Nada Gouda, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range("G38").Value = 'Dr. Abeer Ragab, Dr. Menna tu''Alllah Mohammad'
$ws.Range("G40").Value = 'Dr. Amany Raafat, Dr. Nourhan Osama, Dr. Nahed Mosaad, Dr. Eman M. Abo-Sakaya, Dr. Sara Atawia, Dr. Maryam Ahmad, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Merna Mahrous, Dr. Yasmeena Fattoh, Dr. Merna Said, Dr. Marina Youhanna, Dr. Mai Mustafa'
$ws.Range("G41").Value = 'Dr. Nadia Mostafa, Dr. Amany Raafat, Dr. Nourhan Osama, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Merna Mahrous, Dr. Dina Adel, Dr. Amira Ibrahim'
$ws.Range("G43").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G44").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G47").Value = 'Dr. Aya Alaa-Eldein, Dr. Afaf Abdallah'
$ws.Range("G49").Value = 'Dr. Neveen Nashaat, Dr. Monica, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Naema Gomaa'
$ws.Range("G50").Value = 'Dr. Amira Sobhy, Administrator, Dr. Asmaa Reda, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Nahla Nagiub'
$ws.Range("G51").Value = 'Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Administrator, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Rana Abo-Zaid'
$ws.Range("G52").Value = 'Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G53").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Nada Gouda, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range("G54").Value = 'Dr. Abeer Ragab, Dr. Menna tu''Alllah Mohammad'
$ws.Range("G56").Value = 'Dr. Amany Raafat, Dr. Nourhan Osama, Dr. Nahed Mosaad, Dr. Eman M. Abo-Sakaya, Dr. Sara Atawia, Dr. Maryam Ahmad, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Merna Mahrous, Dr. Yasmeena Fattoh, Dr. Merna Said, Dr. Marina Youhanna, Dr. Mai Mustafa'
$ws.Range("G57").Value = 'Dr. Nadia Mostafa, Dr. Amany Raafat, Dr. Nourhan Osama, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Merna Mahrous, Dr. Dina Adel, Dr. Amira Ibrahim'
$ws.Range("G59").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G60").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G63").Value = 'Dr. Aya Alaa-Eldein, Dr. Afaf Abdallah'
$ws.Range("G65").Value = 'Dr. Neveen Nashaat, Dr. Monica, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Naema Gomaa'
$ws.Range("G66").Value = 'Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub'
$ws.Range("G67").Value = 'Dr. Asmaa Reda, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid'
$ws.Range("G68").Value = 'Dr. Nourhan Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Alshimaa Atef'
$ws.Range("G69").Value = 'Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab'
$ws.Range("G70").Value = 'Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Nada Gouda'
$ws.Range("G71").Value = 'Dr. Nourhan Mohammad, Dr. Sara Nabil'
$ws.Range("G72").Value = 'Dr. Amany Raafat, Dr. Nourhan Osama, Dr. Nahed Mosaad, Dr. Eman M. Abo-Sakaya, Dr. Sara Atawia, Dr. Maryam Ahmad, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Merna Mahrous, Dr. Yasmeena Fattoh, Dr. Merna Said, Dr. Marina Youhanna, Dr. Mai Mustafa'
$ws.Range("G73").Value = 'Dr. Dalia Mohammad Abd Al-Salam, Dr. Nahed Mosaad, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Merna Said'
$ws.Range("G74").Value = 'Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat'
$ws.Range("G75").Value = 'Dr. Mona Ibrahim Hussein, Dr. Alaa Ashraf'
$ws.Range("G76").Value = 'Dr. Mona Ibrahim Hussein, Dr. Alaa Ashraf'
$ws.Range("G79").Value = 'Dr. Walaa Ghanima, Dr. Amr Saeed'
$ws.Range("G80").Value = 'Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Aya Hanafy, Dr. Eman Mohammad Al'
$ws.Range("G81").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G82").Value = 'Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub'
$ws.Range("G83").Value = 'Dr. Mohammad El-Tanany, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad'
$ws.Range("G84").Value = 'Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G85").Value = 'Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab'
$ws.Range("G86").Value = 'Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Nada Gouda'
$ws.Range("G87").Value = 'Dr. Nourhan Mohammad, Dr. Sara Nabil'
$ws.Range("G88").Value = 'Dr. Amany Raafat, Dr. Nourhan Osama, Dr. Nahed Mosaad, Dr. Eman M. Abo-Sakaya, Dr. Sara Atawia, Dr. Maryam Ahmad, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Merna Mahrous, Dr. Yasmeena Fattoh, Dr. Merna Said, Dr. Marina Youhanna, Dr. Mai Mustafa'
$ws.Range("G89").Value = 'Dr. Dalia Mohammad Abd Al-Salam, Dr. Nahed Mosaad, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Merna Said'
$ws.Range("G90").Value = 'Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat'
$ws.Range("G91").Value = 'Dr. Mona Ibrahim Hussein, Dr. Alaa Ashraf'
$ws.Range("G92").Value = 'Dr. Mona Ibrahim Hussein, Dr. Alaa Ashraf'
$ws.Range("G95").Value = 'Dr. Walaa Ghanima, Dr. Amr Saeed'
$ws.Range("G96").Value = 'Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Aya Hanafy, Dr. Eman Mohammad Al'
$ws.Range("G98").Value = 'Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Nahla Nagiub'
$ws.Range("G100").Value = 'Dr. Nourhan Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Alshimaa Atef'
$ws.Range("G101").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Nada Gouda, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range("G102").Value = 'Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Nada Gouda'
$ws.Range("G103").Value = 'Dr. Arwa Elnagar, Dr. Shimaa Ashraf, Dr. Aya Saeed'
$ws.Range("G104").Value = 'Dr. Nourhan Osama, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim'
$ws.Range("G105").Value = 'Dr. Nadia Mostafa, Dr. Amany Raafat, Dr. Nourhan Osama, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Merna Mahrous, Dr. Dina Adel, Dr. Amira Ibrahim'
$ws.Range("G107").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G108").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G112").Value = 'Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Nahla, Dr. Remon, Dr. Youstina Magdy, Dr. Yassmen Ahmad'
$ws.Range("G114").Value = 'Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Nahla Nagiub'
$ws.Range("G116").Value = 'Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Manar Montaser, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Majorelle Magdy'
$ws.Range("G117").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Nada Gouda, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range("G118").Value = 'Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Nada Gouda'
$ws.Range("G119").Value = 'Dr. Arwa Elnagar, Dr. Shimaa Ashraf, Dr. Aya Saeed'
$ws.Range("G120").Value = 'Dr. Nourhan Osama, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim'
$ws.Range("G121").Value = 'Dr. Nadia Mostafa, Dr. Amany Raafat, Dr. Nourhan Osama, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Merna Mahrous, Dr. Dina Adel, Dr. Amira Ibrahim'
$ws.Range("G123").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G124").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G125").Value = 'Dr. Walaa Ghanima, Dr. Nancy Abd Al-Shafy'
$ws.Range("G128").Value = 'Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Nahla, Dr. Remon, Dr. Youstina Magdy, Dr. Yassmen Ahmad'
